# [MOSIP-14336] Updating Masterdata Utility
# Rebuild the "language" masterdata upload sheet with the new column layout
# (leading numeric index column + renamed/relocated headers) and seed rows
# for English, Arabic and French.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 used to hold the "code" header (with the bold/bordered header format);
# that header now lives one column over (B1), so reset A1 back to the
# workbook's default "Normal" style and blank it out entirely.
$a1 = $ws.Range("A1")
$a1.Style = "Normal"
$a1.Value = $null

# --- Header row (row 1) ---------------------------------------------------
# Headers now start at column B; "nativeName"/"isActive" become
# "native_name"/"is_active", shifted one column to the right, and column A
# is reserved for a new leading numeric index column.
$ws.Cells.Item(1, 2).Value = "code"
$ws.Cells.Item(1, 3).Value = "name"
$ws.Cells.Item(1, 4).Value = "family"
$ws.Cells.Item(1, 5).Value = "native_name"
$ws.Cells.Item(1, 6).Value = "is_active"

# Re-apply the bold / centered-top / thin-bordered header formatting (the
# same formatting the original A1:E1 header row carried) to the relocated
# header cells B1:E1 (all within the sheet's original used range).
$srcRange = $ws.Range("B1:E1")
$srcRange.Borders.LineStyle = 1
$srcRange.HorizontalAlignment = -4108
$srcRange.VerticalAlignment = -4160
$srcRange.Font.Bold = $true

# F1 is a brand new cell outside the original used range, so copy the
# header formatting over (format-only paste) rather than rebuild it from
# individual properties, to keep reusing the very same header cell style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# --- Data rows --------------------------------------------------------------
# index, code, name, family, native_name, is_active  -- the new leading
# numeric index column carries the same header formatting as row 1.
$data = @(
    @(0, "eng", "English", "Indo-European", "English", $true),
    @(1, "ara", "Arabic", "Afro-Asiatic", "العَرَبِيَّة‎", $true),
    @(2, "fra", "French", "Indo-European", "français", $true)
)

$r = 2
foreach ($row in $data) {
    $idxCell = $ws.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $ws.Range("B1").Copy() | Out-Null
    $idxCell.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$excel.CutCopyMode = $false
